# Daily attendance processing - 2026-01-18 12:51:51
# Swap the order of "dnasr281@gmail.com" and "System" in column G
# (Recorded By) wherever both appear together as
# "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldText) {
        $cell.Value = $newText
    }
}
